$wb = $excel.ActiveWorkbook

# --- Sheet "Prix Spot": insert a new "08-dec" column before the "01-oct." column (EK) ---
$ws1 = $wb.Worksheets.Item("Prix Spot")

# Insert a new blank column at EK, shifting EK:FO (and everything after) one column to the right.
$ws1.Columns("EK:EK").Insert()

# Header cell for the newly inserted column.
$ws1.Range("EK1").Value = "08-dec"

# Data cells for the newly inserted column (rows 2-25) are all "-" (no data for that day yet).
for ($r = 2; $r -le 25; $r++) {
    $ws1.Cells.Item($r, 141).Value = "-"
}

# --- Sheet "Gaz": append two new daily rows ---
# (dates are kept as plain text, like the rest of column A, via a quote-prefixed
#  value followed by resetting the style so Excel doesn't auto-convert/format them
#  as a date serial number)
$wsGaz = $wb.Worksheets.Item("Gaz")
$wsGaz.Range("A171").Value = "'2025-12-06"
$wsGaz.Range("A171").Style = "Normal"
$wsGaz.Range("B171").Value = 25.905
$wsGaz.Range("A172").Value = "'2025-12-07"
$wsGaz.Range("A172").Style = "Normal"
$wsGaz.Range("B172").Value = 25.905

# --- Sheet "CO2": append two new daily rows ---
$wsCo2 = $wb.Worksheets.Item("CO2")
$wsCo2.Range("A171").Value = "'2025-12-06"
$wsCo2.Range("A171").Style = "Normal"
$wsCo2.Range("B171").Value = 81.78
$wsCo2.Range("A172").Value = "'2025-12-07"
$wsCo2.Range("A172").Style = "Normal"
$wsCo2.Range("B172").Value = 81.78
